# Edit script: adds a new "PO Forecast" worksheet with forecast data,
# and renames the "Requested quantity" header on the two existing
# sheets to their new, more specific names.

$wb = $excel.ActiveWorkbook

# --- 1) Rename existing headers -------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the new "PO Forecast" sheet after the existing sheets ---------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "PO Forecast"

# Header row
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$newSheet.Cells.Item(2, 1).Value = 45025.99999999999
$newSheet.Cells.Item(2, 2).Value = 45
$newSheet.Cells.Item(2, 3).Value = -135.459919342918
$newSheet.Cells.Item(2, 4).Value = 231.8710996451092
$newSheet.Cells.Item(3, 1).Value = 45032.99999999999
$newSheet.Cells.Item(3, 2).Value = 47
$newSheet.Cells.Item(3, 3).Value = -125.9214927052753
$newSheet.Cells.Item(3, 4).Value = 221.6824527897937
$newSheet.Cells.Item(4, 1).Value = 45046.99999999999
$newSheet.Cells.Item(4, 2).Value = 52
$newSheet.Cells.Item(4, 3).Value = -129.5316772482637
$newSheet.Cells.Item(4, 4).Value = 240.3259582548227
$newSheet.Cells.Item(5, 1).Value = 45060.99999999999
$newSheet.Cells.Item(5, 2).Value = 56
$newSheet.Cells.Item(5, 3).Value = -127.8338262798837
$newSheet.Cells.Item(5, 4).Value = 240.5259328454014
$newSheet.Cells.Item(6, 1).Value = 45074.99999999999
$newSheet.Cells.Item(6, 2).Value = 61
$newSheet.Cells.Item(6, 3).Value = -134.311135118172
$newSheet.Cells.Item(6, 4).Value = 241.8977426201488
$newSheet.Cells.Item(7, 1).Value = 45081.99999999999
$newSheet.Cells.Item(7, 2).Value = 63
$newSheet.Cells.Item(7, 3).Value = -121.268380924607
$newSheet.Cells.Item(7, 4).Value = 239.7467713340835
$newSheet.Cells.Item(8, 1).Value = 45088.99999999999
$newSheet.Cells.Item(8, 2).Value = 65
$newSheet.Cells.Item(8, 3).Value = -122.6061690879
$newSheet.Cells.Item(8, 4).Value = 237.6025393631747
$newSheet.Cells.Item(9, 1).Value = 45102.99999999999
$newSheet.Cells.Item(9, 2).Value = 70
$newSheet.Cells.Item(9, 3).Value = -114.3557824932196
$newSheet.Cells.Item(9, 4).Value = 257.2798687601219
$newSheet.Cells.Item(10, 1).Value = 45109.99999999999
$newSheet.Cells.Item(10, 2).Value = 72
$newSheet.Cells.Item(10, 3).Value = -101.2968504248589
$newSheet.Cells.Item(10, 4).Value = 258.3974195730356
$newSheet.Cells.Item(11, 1).Value = 45116.99999999999
$newSheet.Cells.Item(11, 2).Value = 74
$newSheet.Cells.Item(11, 3).Value = -127.1412877056628
$newSheet.Cells.Item(11, 4).Value = 268.0391426632834
$newSheet.Cells.Item(12, 1).Value = 45123.99999999999
$newSheet.Cells.Item(12, 2).Value = 77
$newSheet.Cells.Item(12, 3).Value = -126.2751649772667
$newSheet.Cells.Item(12, 4).Value = 250.3018014095274
$newSheet.Cells.Item(13, 1).Value = 45130.99999999999
$newSheet.Cells.Item(13, 2).Value = 79
$newSheet.Cells.Item(13, 3).Value = -116.1110010025477
$newSheet.Cells.Item(13, 4).Value = 242.8287405462692
$newSheet.Cells.Item(14, 1).Value = 45137.99999999999
$newSheet.Cells.Item(14, 2).Value = 81
$newSheet.Cells.Item(14, 3).Value = -100.846505954682
$newSheet.Cells.Item(14, 4).Value = 253.5430986377548
$newSheet.Cells.Item(15, 1).Value = 45165.99999999999
$newSheet.Cells.Item(15, 2).Value = 90
$newSheet.Cells.Item(15, 3).Value = -118.6422824123685
$newSheet.Cells.Item(15, 4).Value = 275.6864566502211
$newSheet.Cells.Item(16, 1).Value = 45179.99999999999
$newSheet.Cells.Item(16, 2).Value = 95
$newSheet.Cells.Item(16, 3).Value = -85.89644026904887
$newSheet.Cells.Item(16, 4).Value = 284.1068497560404
$newSheet.Cells.Item(17, 1).Value = 45186.99999999999
$newSheet.Cells.Item(17, 2).Value = 97
$newSheet.Cells.Item(17, 3).Value = -85.14811629552784
$newSheet.Cells.Item(17, 4).Value = 277.9437276270332
$newSheet.Cells.Item(18, 1).Value = 45200.99999999999
$newSheet.Cells.Item(18, 2).Value = 102
$newSheet.Cells.Item(18, 3).Value = -82.70847603334788
$newSheet.Cells.Item(18, 4).Value = 273.0713899038622
$newSheet.Cells.Item(19, 1).Value = 45207.99999999999
$newSheet.Cells.Item(19, 2).Value = 104
$newSheet.Cells.Item(19, 3).Value = -79.47259651133655
$newSheet.Cells.Item(19, 4).Value = 282.7283174725093
$newSheet.Cells.Item(20, 1).Value = 45214.99999999999
$newSheet.Cells.Item(20, 2).Value = 106
$newSheet.Cells.Item(20, 3).Value = -77.63043040066152
$newSheet.Cells.Item(20, 4).Value = 291.1155387786626
$newSheet.Cells.Item(21, 1).Value = 45221.99999999999
$newSheet.Cells.Item(21, 2).Value = 109
$newSheet.Cells.Item(21, 3).Value = -81.62847886840376
$newSheet.Cells.Item(21, 4).Value = 287.0160914274454
$newSheet.Cells.Item(22, 1).Value = 45228.99999999999
$newSheet.Cells.Item(22, 2).Value = 111
$newSheet.Cells.Item(22, 3).Value = -62.38726841297921
$newSheet.Cells.Item(22, 4).Value = 275.7667627368092
$newSheet.Cells.Item(23, 1).Value = 45235.99999999999
$newSheet.Cells.Item(23, 2).Value = 113
$newSheet.Cells.Item(23, 3).Value = -87.69518878495028
$newSheet.Cells.Item(23, 4).Value = 296.1133490654085
$newSheet.Cells.Item(24, 1).Value = 45242.99999999999
$newSheet.Cells.Item(24, 2).Value = 116
$newSheet.Cells.Item(24, 3).Value = -63.78388459256309
$newSheet.Cells.Item(24, 4).Value = 296.2153470113641
$newSheet.Cells.Item(25, 1).Value = 45249.99999999999
$newSheet.Cells.Item(25, 2).Value = 118
$newSheet.Cells.Item(25, 3).Value = -43.31493122723438
$newSheet.Cells.Item(25, 4).Value = 307.7192532528065
$newSheet.Cells.Item(26, 1).Value = 45263.99999999999
$newSheet.Cells.Item(26, 2).Value = 122
$newSheet.Cells.Item(26, 3).Value = -50.1457913452662
$newSheet.Cells.Item(26, 4).Value = 302.7507561763107
$newSheet.Cells.Item(27, 1).Value = 45305.99999999999
$newSheet.Cells.Item(27, 2).Value = 136
$newSheet.Cells.Item(27, 3).Value = -64.49402312040519
$newSheet.Cells.Item(27, 4).Value = 319.8817105694295
$newSheet.Cells.Item(28, 1).Value = 45319.99999999999
$newSheet.Cells.Item(28, 2).Value = 141
$newSheet.Cells.Item(28, 3).Value = -55.6790782136946
$newSheet.Cells.Item(28, 4).Value = 324.5702847539085
$newSheet.Cells.Item(29, 1).Value = 45333.99999999999
$newSheet.Cells.Item(29, 2).Value = 145
$newSheet.Cells.Item(29, 3).Value = -36.69252894755095
$newSheet.Cells.Item(29, 4).Value = 337.8261652636446
$newSheet.Cells.Item(30, 1).Value = 45361.99999999999
$newSheet.Cells.Item(30, 2).Value = 155
$newSheet.Cells.Item(30, 3).Value = -17.46204079211316
$newSheet.Cells.Item(30, 4).Value = 332.8422782904261
$newSheet.Cells.Item(31, 1).Value = 45368.99999999999
$newSheet.Cells.Item(31, 2).Value = 157
$newSheet.Cells.Item(31, 3).Value = -29.37424080809727
$newSheet.Cells.Item(31, 4).Value = 340.3670621113837
$newSheet.Cells.Item(32, 1).Value = 45375.99999999999
$newSheet.Cells.Item(32, 2).Value = 159
$newSheet.Cells.Item(32, 3).Value = -27.88555923886067
$newSheet.Cells.Item(32, 4).Value = 358.5358125451778
$newSheet.Cells.Item(33, 1).Value = 45403.99999999999
$newSheet.Cells.Item(33, 2).Value = 168
$newSheet.Cells.Item(33, 3).Value = -20.21478101579965
$newSheet.Cells.Item(33, 4).Value = 352.2751035483438
$newSheet.Cells.Item(34, 1).Value = 45410.99999999999
$newSheet.Cells.Item(34, 2).Value = 171
$newSheet.Cells.Item(34, 3).Value = -7.269797726477624
$newSheet.Cells.Item(34, 4).Value = 359.7910509719134
$newSheet.Cells.Item(35, 1).Value = 45445.99999999999
$newSheet.Cells.Item(35, 2).Value = 182
$newSheet.Cells.Item(35, 3).Value = -11.69845039633446
$newSheet.Cells.Item(35, 4).Value = 351.1005866498431
$newSheet.Cells.Item(36, 1).Value = 45466.99999999999
$newSheet.Cells.Item(36, 2).Value = 189
$newSheet.Cells.Item(36, 3).Value = 4.516801204204577
$newSheet.Cells.Item(36, 4).Value = 373.1531852335118
$newSheet.Cells.Item(37, 1).Value = 45473.99999999999
$newSheet.Cells.Item(37, 2).Value = 191
$newSheet.Cells.Item(37, 3).Value = 10.4021572689153
$newSheet.Cells.Item(37, 4).Value = 394.0609497390744
$newSheet.Cells.Item(38, 1).Value = 45487.99999999999
$newSheet.Cells.Item(38, 2).Value = 196
$newSheet.Cells.Item(38, 3).Value = 8.895091648602632
$newSheet.Cells.Item(38, 4).Value = 391.6362146242161
$newSheet.Cells.Item(39, 1).Value = 45494.99999999999
$newSheet.Cells.Item(39, 2).Value = 198
$newSheet.Cells.Item(39, 3).Value = 13.8668573788341
$newSheet.Cells.Item(39, 4).Value = 372.1774957150297
$newSheet.Cells.Item(40, 1).Value = 45501.99999999999
$newSheet.Cells.Item(40, 2).Value = 200
$newSheet.Cells.Item(40, 3).Value = 21.75983861629605
$newSheet.Cells.Item(40, 4).Value = 381.870923583422
$newSheet.Cells.Item(41, 1).Value = 45508.99999999999
$newSheet.Cells.Item(41, 2).Value = 203
$newSheet.Cells.Item(41, 3).Value = 17.23647147462314
$newSheet.Cells.Item(41, 4).Value = 391.9807246602838
$newSheet.Cells.Item(42, 1).Value = 45515.99999999999
$newSheet.Cells.Item(42, 2).Value = 205
$newSheet.Cells.Item(42, 3).Value = 21.12199312126419
$newSheet.Cells.Item(42, 4).Value = 380.1202825251706
$newSheet.Cells.Item(43, 1).Value = 45522.99999999999
$newSheet.Cells.Item(43, 2).Value = 207
$newSheet.Cells.Item(43, 3).Value = 37.11113904911878
$newSheet.Cells.Item(43, 4).Value = 399.1320433573617
$newSheet.Cells.Item(44, 1).Value = 45529.99999999999
$newSheet.Cells.Item(44, 2).Value = 209
$newSheet.Cells.Item(44, 3).Value = 32.1826062010712
$newSheet.Cells.Item(44, 4).Value = 403.4802784587497
$newSheet.Cells.Item(45, 1).Value = 45536.99999999999
$newSheet.Cells.Item(45, 2).Value = 212
$newSheet.Cells.Item(45, 3).Value = 27.51218060227043
$newSheet.Cells.Item(45, 4).Value = 407.0433640317528
$newSheet.Cells.Item(46, 1).Value = 45543.99999999999
$newSheet.Cells.Item(46, 2).Value = 214
$newSheet.Cells.Item(46, 3).Value = 21.17189894597619
$newSheet.Cells.Item(46, 4).Value = 410.5244254475149
$newSheet.Cells.Item(47, 1).Value = 45557.99999999999
$newSheet.Cells.Item(47, 2).Value = 219
$newSheet.Cells.Item(47, 3).Value = 35.08338831309972
$newSheet.Cells.Item(47, 4).Value = 390.9804223000712
$newSheet.Cells.Item(48, 1).Value = 45564.99999999999
$newSheet.Cells.Item(48, 2).Value = 221
$newSheet.Cells.Item(48, 3).Value = 36.37290701904812
$newSheet.Cells.Item(48, 4).Value = 401.2186025766567
$newSheet.Cells.Item(49, 1).Value = 45571.99999999999
$newSheet.Cells.Item(49, 2).Value = 223
$newSheet.Cells.Item(49, 3).Value = 25.16807166924611
$newSheet.Cells.Item(49, 4).Value = 408.4521207543486
$newSheet.Cells.Item(50, 1).Value = 45578.99999999999
$newSheet.Cells.Item(50, 2).Value = 225
$newSheet.Cells.Item(50, 3).Value = 38.94229926762693
$newSheet.Cells.Item(50, 4).Value = 409.7753329397426
$newSheet.Cells.Item(51, 1).Value = 45585.99999999999
$newSheet.Cells.Item(51, 2).Value = 228
$newSheet.Cells.Item(51, 3).Value = 37.839983753106
$newSheet.Cells.Item(51, 4).Value = 408.6371156899102
$newSheet.Cells.Item(52, 1).Value = 45606.99999999999
$newSheet.Cells.Item(52, 2).Value = 235
$newSheet.Cells.Item(52, 3).Value = 51.81784183557672
$newSheet.Cells.Item(52, 4).Value = 419.46053010466
$newSheet.Cells.Item(53, 1).Value = 45613.99999999999
$newSheet.Cells.Item(53, 2).Value = 237
$newSheet.Cells.Item(53, 3).Value = 57.06987161347275
$newSheet.Cells.Item(53, 4).Value = 421.6788362494598
$newSheet.Cells.Item(54, 1).Value = 45620.99999999999
$newSheet.Cells.Item(54, 2).Value = 239
$newSheet.Cells.Item(54, 3).Value = 68.79096578692035
$newSheet.Cells.Item(54, 4).Value = 419.3427462058393
$newSheet.Cells.Item(55, 1).Value = 45627.99999999999
$newSheet.Cells.Item(55, 2).Value = 241
$newSheet.Cells.Item(55, 3).Value = 70.98936151175647
$newSheet.Cells.Item(55, 4).Value = 421.617934822798
$newSheet.Cells.Item(56, 1).Value = 45634.99999999999
$newSheet.Cells.Item(56, 2).Value = 244
$newSheet.Cells.Item(56, 3).Value = 55.36654740858081
$newSheet.Cells.Item(56, 4).Value = 435.3051513278211
$newSheet.Cells.Item(57, 1).Value = 45641.99999999999
$newSheet.Cells.Item(57, 2).Value = 246
$newSheet.Cells.Item(57, 3).Value = 54.05602942797178
$newSheet.Cells.Item(57, 4).Value = 423.5533209189615
$newSheet.Cells.Item(58, 1).Value = 45648.99999999999
$newSheet.Cells.Item(58, 2).Value = 248
$newSheet.Cells.Item(58, 3).Value = 64.32345041327392
$newSheet.Cells.Item(58, 4).Value = 430.777103711612
$newSheet.Cells.Item(59, 1).Value = 45655.99999999999
$newSheet.Cells.Item(59, 2).Value = 251
$newSheet.Cells.Item(59, 3).Value = 71.51871812745064
$newSheet.Cells.Item(59, 4).Value = 434.3208700643095
$newSheet.Cells.Item(60, 1).Value = 45662.99999999999
$newSheet.Cells.Item(60, 2).Value = 253
$newSheet.Cells.Item(60, 3).Value = 83.18768443234821
$newSheet.Cells.Item(60, 4).Value = 444.3784683710468

# --- 3) Match formatting used on the other sheets --------------------------
# Header style (bold, centered, bordered) -- copy from an existing header cell
$wsWeekly.Range("A1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Date number-format style for the "ds" column -- copy from an existing date cell
$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A60").PasteSpecial(-4122)

Write-Output "PO Forecast sheet added and headers renamed."
